$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from LOINC")

# Current layout (rows 1-13):
#   row 1   : header  (Concept / Description)
#   rows 2-11: LOINC codes, column B blank
#   row 12  : blank separator row (both columns blank)
#   row 13  : System URI / http://loinc.org footer
#
# Target layout (rows 1-14): a new LOINC concept row "42348-3" is added
# right after the existing code rows (as row 12); the blank separator row
# and the footer row are pushed down to rows 13 and 14 respectively.

# Remember the rows that need to move down before we overwrite anything.
$oldA12 = $ws.Cells.Item(12, 1).Value()
$oldB12 = $ws.Cells.Item(12, 2).Value()
$oldA13 = $ws.Cells.Item(13, 1).Value()
$oldB13 = $ws.Cells.Item(13, 2).Value()

# Footer row (old row 13) -> row 14.
$ws.Cells.Item(14, 1).Value = $oldA13
$ws.Cells.Item(14, 2).Value = $oldB13

# Blank separator row (old row 12) -> row 13.
$ws.Cells.Item(13, 1).Value = $oldA12
$ws.Cells.Item(13, 2).Value = $oldB12

# New concept row in row 12: LOINC code only, description left blank.
$ws.Cells.Item(12, 1).Value = "42348-3"
$ws.Cells.Item(12, 2).Value = ""

# Give the newly-populated row 14 the same formatting as the other data
# rows (rows 12/13 already inherited it). Use a formats-only paste so we
# don't create a duplicate style entry.
$ws.Range("A11:B11").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
